{"js": "// Apply the review-434 rewrite: update the title/date line, rewrite the\n// body paragraphs, and append four new paragraphs (incl. the new source link)\n// that extend the review, replacing the single old arxiv link paragraph.\n\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// Paragraph 1 - date/title line (kept as two runs joined by a manual line break)\nparas.items[0].insertText(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 05.04.25\\u000bGIVT: Generative Infinite-Vocabulary Transformers\", Word.InsertLocation.replace);\n\n// Paragraphs 2-5 - rewritten review body\nparas.items[1].insertText(\"\u05d4\u05d9\u05d5\u05dd \u05d7\u05d5\u05d6\u05e8\u05d9\u05dd \u05db\u05de\u05d4 \u05e9\u05e0\u05d9\u05dd \u05e9\u05e0\u05d9\u05dd \u05d0\u05d7\u05d5\u05e8\u05d4 \u05d1\u05d5 \u05de\u05d9\u05dc\u05d9\u05dd VAE, VQ-VAE, VQ-GAN \u05d4\u05d9\u05d5 \u05de\u05d5\u05e9\u05db\u05d9\u05dd \u05d0\u05d5\u05ea\u05d4 \u05ea\u05e9\u05d5\u05de\u05ea \u05dc\u05d1 \u05db\u05de\u05d5 \u05e9\u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05d4\u05d9\u05d5\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd (\u05d0\u05de\u05e0\u05dd \u05e4\u05d7\u05d5\u05ea \u05de\u05d0\u05d2'\u05e0\u05d8\u05d9\u05dd \u05d0\u05d1\u05dc \u05d1\u05db\u05dc \u05d6\u05d0\u05ea).  \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05de\u05e6\u05d9\u05e2 \u05e9\u05db\u05dc\u05d5\u05dc \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05dc-VQ-VAE \u05e9\u05de\u05e9\u05da \u05d0\u05ea \u05e2\u05d9\u05e0\u05d9\u05d9 \u05db\u05d9 \u05db\u05d0\u05de\u05d5\u05e8 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05d1\u05e0\u05d5\u05e9\u05d0 \u05d6\u05d4 \u05d4\u05e4\u05db\u05d5 \u05dc\u05d4\u05d9\u05d5\u05ea \u05f4\u05e6\u05d9\u05e4\u05d5\u05e8 \u05e0\u05d3\u05d9\u05e8\u05d4\u05f4 \u05d1\u05e0\u05d5\u05e3 \u05e9\u05dc\u05e0\u05d5 (\u05e9\u05dc AI).\", Word.InsertLocation.replace);\nparas.items[2].insertText(\"\u05e7\u05d5\u05d3\u05dd \u05db\u05dc \u05d0\u05ea\u05df \u05d4\u05e7\u05d3\u05de\u05d4 \u05e7\u05e6\u05e8 \u05dc\u05d2\u05d1\u05d9 VQ-VAE. \u05e0\u05ea\u05d7\u05d9\u05dc \u05d0\u05ea \u05d4\u05d4\u05e1\u05d1\u05e8 \u05de-VAE \u05e9\u05d6\u05d4 \u05e8\u05d0\u05e9\u05d9 \u05ea\u05d9\u05d1\u05d5\u05ea \u05e9\u05dc Variational AutoEncoder \u05e9\u05d4\u05d5\u05de\u05e6\u05d0 \u05d0\u05d9 \u05e9\u05dd \u05d1-2014 \u05e2\u05dc \u05d9\u05d3\u05d9 Kingma \u05d4\u05d0\u05d2\u05d3\u05d9. \u05d1\u05d2\u05d3\u05d5\u05dc VAE \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05e9\u05ea\u05d9 \u05e8\u05e9\u05ea\u05d5\u05ea, \u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d3\u05e7\u05d5\u05d3\u05e8 \u05db\u05d0\u05e9\u05e8 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05de\u05e4\u05d9\u05e7 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 (\u05d0\u05d5 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2) \u05e9\u05dc \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05db\u05d0\u05e9\u05e8 \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05dc\u05ea\u05de\u05d5\u05e0\u05d4. \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d4\u05d5\u05d0 \u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05dc \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d4\u05d2\u05d0\u05d5\u05e1\u05d9\u05ea (\u05d5\u05e7\u05d8\u05d5\u05e8 \u05ea\u05d5\u05d7\u05dc\u05d5\u05ea \u05d5\u05de\u05d8\u05e8\u05d9\u05e6\u05ea \u05e7\u05d5\u05d5\u05e8\u05d9\u05d0\u05e0\u05e1 \u05d0\u05dc\u05db\u05e1\u05d5\u05e0\u05d9\u05ea) \u05de\u05de\u05e0\u05d4 \u05d3\u05d5\u05d2\u05de\u05d9\u05dd \u05d0\u05ea \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05d4\u05de\u05d5\u05d6\u05df \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05dc\u05e9\u05d7\u05d6\u05d5\u05e8 \u05ea\u05de\u05d5\u05e0\u05ea \u05d4\u05e7\u05dc\u05d8 \u05dc\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8. \", Word.InsertLocation.replace);\nparas.items[3].insertText(\"\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05dc\u05d5\u05e1 \u05e9\u05dc VAE \u05e0\u05d1\u05e0\u05d9\u05ea \u05e2\u05dc \u05d1\u05e1\u05d9\u05e1 ELBO (\u05e9\u05d6\u05d4 Evidence Lower Bound) \u05d5\u05de\u05db\u05d9\u05dc\u05d4 2 \u05d0\u05d9\u05d1\u05e8\u05d9\u05dd. \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05dc\u05d5\u05e1 \u05d4\u05e9\u05d7\u05d6\u05d5\u05e8 \u05e9\u05dc-VAE \u05e7\u05dc\u05d0\u05e1\u05d9 \u05d4\u05d5\u05d0 \u05d4\u05e0\u05d5\u05e8\u05de\u05d4 \u05e9\u05dc \u05d4\u05e4\u05e8\u05e9 \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d4\u05de\u05e9\u05d5\u05d7\u05d6\u05e8\u05ea \u05d1\u05d9\u05d7\u05e1 \u05dc\u05ea\u05de\u05d5\u05e0\u05d4 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9\u05ea (\u05d1\u05d2\u05e8\u05e1\u05d0\u05d5\u05ea \u05de\u05ea\u05e7\u05d3\u05de\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d4\u05ea\u05d5\u05d5\u05e1\u05e4\u05d5 \u05dc\u05d6\u05d4 \u05dc\u05d5\u05e1 perceptual \u05d5\u05dc\u05d5\u05e1 \u05d1\u05e1\u05d2\u05e0\u05d5\u05df GAN) \u05d5\u05d4\u05d0\u05d9\u05d1\u05e8 \u05d4\u05e9\u05e0\u05d9 \u05d4\u05d5\u05d0 KL divergence \u05d1\u05d9\u05df \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e9\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05d4\u05de\u05d5\u05e4\u05e7 \u05de\u05d4\u05d3\u05d0\u05d8\u05d4 (\u05d4\u05de\u05d9\u05d5\u05e6\u05d2 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d5\u05e7\u05d8\u05d5\u05e8 \u05ea\u05d5\u05d7\u05dc\u05d5\u05ea \u05d5\u05de\u05d8\u05e8\u05d9\u05e6\u05ea \u05e7\u05d5\u05d5\u05e8\u05d9\u05d0\u05e0\u05e1 \u05d0\u05dc\u05db\u05e1\u05d5\u05e0\u05d9\u05ea) \u05dc\u05d1\u05d9\u05df \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e0\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea. \u05d1\u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1 \u05d0\u05e0\u05d5 \u05d3\u05d5\u05d2\u05de\u05d9\u05dd \u05d5\u05e7\u05d8\u05d5\u05e8 \u05dc\u05d8\u05e0\u05d8\u05d9 \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e0\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea \u05d5\u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05d5\u05ea\u05d5 \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8.\", Word.InsertLocation.replace);\nparas.items[4].insertText(\"\u05e9\u05db\u05dc\u05d5\u05dc \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e9\u05d4\u05e4\u05da \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d0\u05d5\u05d3 \u05e4\u05d5\u05e4\u05d5\u05dc\u05e8\u05d9 \u05e9\u05dc VAE \u05d4\u05d5\u05d0 VQ-VAE. \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05d4\u05d2\u05d3\u05d9\u05e8 \u05de\u05e8\u05d7\u05d1 \u05dc\u05d8\u05e0\u05d8\u05d9 \u05d1\u05ea\u05d5\u05e8 \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea - \u05d0\u05dc\u05d0 \u05de\u05d2\u05d3\u05d9\u05e8 \u05d0\u05d5\u05ea\u05d5 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d9\u05e1\u05e7\u05e8\u05d8\u05d9\u05ea. \u05db\u05dc \u05e4\u05d0\u05e5' \u05d1\u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05ea\u05d5\u05d0\u05e8(\u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9) \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d5\u05e7\u05d8\u05d5\u05e8 \u05de\u05d4-codebook \u05d1\u05d2\u05d5\u05d3\u05dc \u05e1\u05d5\u05e4\u05d9 \u05e9\u05de\u05d0\u05d5\u05de\u05df \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d4\u05d3\u05e7\u05d5\u05d3\u05e8. \u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05e9 \u05de\u05e1\u05e4\u05e8 \u05e1\u05d5\u05e4\u05d9 \u05e9\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05e4\u05d0\u05e5' (\u05e0\u05d7\u05d6\u05d5\u05e8 \u05e2\u05dc \u05d6\u05d4 \u05e2\u05d5\u05d3 \u05de\u05e2\u05d8). \u05db\u05d0\u05e9\u05e8 \u05d0\u05d9\u05de\u05d5\u05df \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8, \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05de\u05d4-codebook \u05de\u05e1\u05ea\u05d9\u05d9\u05dd \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05e0\u05d5\u05e1\u05e3 \u05dc\u05d7\u05d9\u05d6\u05d5\u05d9 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc \u05e4\u05d0\u05e6'\u05d9\u05dd, \u05e2\u05dc \u05db\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8. \u05de\u05d5\u05d3\u05dc \u05d6\u05d4 (\u05e0\u05d2\u05d9\u05d3 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8) \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05d7\u05d6\u05d5\u05ea \u05d1\u05d0\u05d5\u05e4\u05df \u05d0\u05d5\u05d8\u05d5\u05e8\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9 \u05d0\u05ea \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05de\u05d4-codebook (\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05e1\u05e4\u05e8\u05d5) \u05e9\u05dc \u05d4\u05e4\u05d0\u05e5' \u05d4\u05d1\u05d0 \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05d4\u05e7\u05d5\u05d3\u05de\u05d9\u05dd \u05e9\u05db\u05d1\u05e8 \u05d2\u05d5\u05e0\u05e8\u05d8\u05d5. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05d4\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05dc\u05d2\u05e0\u05e8\u05d5\u05d8 \u05d3\u05d0\u05d8\u05d4 (\u05ea\u05de\u05d5\u05e0\u05d4).\", Word.InsertLocation.replace);\n\n// Paragraph 6 - was the old arxiv link, now continues the rewritten review text\nparas.items[5].insertText(\"\u05db\u05d0\u05de\u05d5\u05e8 \u05d9\u05e9 \u05de\u05e1\u05e4\u05e8 \u05de\u05d5\u05d2\u05d1\u05dc \u05d5\u05e1\u05d5\u05e4\u05d9 \u05e9\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05e4\u05d0\u05e5' \u05d5\u05e9\u05d6\u05d4 \u05d3\u05d9 \u05de\u05d2\u05d1\u05d9\u05dc \u05d0\u05ea \u05d4\u05e2\u05d5\u05e9\u05e8 \u05d4\u05e1\u05de\u05e0\u05d8\u05d9 \u05e9\u05dc \u05d4\u05ea\u05de\u05d5\u05e0\u05d5\u05ea \u05e9-VQ-VAE \u05d5\u05e9\u05d9\u05d8\u05d5\u05ea \u05d3\u05d5\u05de\u05d5\u05ea \u05de\u05e1\u05d5\u05d2\u05dc\u05d5\u05ea \u05dc\u05d2\u05e0\u05e8\u05d8. \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05d4\u05de\u05e7\u05d5\u05dd \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05e2\u05db\u05e9\u05d9\u05d5 \u05de\u05d7\u05d3\u05e9 - \u05d4\u05d5\u05d0 \u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05e2\u05d1\u05d5\u05e8 \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e8\u05e6\u05d9\u05e3 (\u05dc\u05dc\u05d0 codebook) \u05e9\u05dc \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd. \u05d0\u05d1\u05dc \u05d0\u05d9\u05da \u05d0\u05e4\u05e9\u05e8 \u05dc\u05e2\u05e9\u05d5\u05ea \u05d6\u05d0\u05ea? \u05e0\u05d6\u05db\u05d9\u05e8 \u05e9-VQ-VAE \u05d0\u05e0\u05d5 \u05db\u05dc \u05e4\u05e2\u05dd \u05d7\u05d5\u05d6\u05d9\u05dd \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d0\u05dc\u05d9\u05ea \u05de\u05e2\u05dc \u05d4-codebook \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05d1\u05de\u05d5\u05d3\u05dc \u05d0\u05d5\u05d8\u05d5\u05e8\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9 \u05d4\u05d9\u05d0 \u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05d1\u05d2\u05d5\u05d3\u05dc \u05e9\u05dc \u05d4-codebook.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Append four new paragraphs after the (rewritten) 6th paragraph, chaining off\n// the paragraph object each `insertParagraph` call returns, ending with the\n// new arxiv source link as the document's final paragraph.\nlet tail = paras.items[5];\ntail = tail.insertParagraph(\" \u05d4\u05d0\u05dd \u05e0\u05d9\u05ea\u05df \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05e9\u05d9\u05d5\u05e6\u05e8 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e8\u05e6\u05d9\u05e4\u05d9\u05dd \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d0\u05d5\u05d8\u05d5\u05e8\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9\u05ea? \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05d9\u05d0 \u05db\u05df - \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05de\u05df \u05d0\u05ea \u05d4-VAE \u05d4\u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9 \u05e9\u05d4\u05e1\u05d1\u05e8\u05ea\u05d9 \u05e2\u05dc\u05d9\u05d5 \u05d1\u05ea\u05d7\u05d9\u05dc\u05ea \u05d4\u05e1\u05e7\u05d9\u05e8\u05d4. \u05d4\u05d3\u05d1\u05e8 \u05d4\u05d6\u05d4 \u05e0\u05e2\u05e9\u05d4 \u05d1\u05e8\u05de\u05d4 \u05e9\u05dc \u05e4\u05d0\u05e5' \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05d4\u05e4\u05d0\u05e6'\u05d9\u05dd \u05e9\u05dc\u05d4. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e9\u05e0\u05d9 \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05e1\u05d9\u05d1\u05ea\u05d9 \u05e9\u05d7\u05d5\u05d6\u05d4 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05d4\u05d8\u05d5\u05e7\u05df \u05d4\u05d1\u05d0 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d7\u05d9\u05d6\u05d5\u05d9 \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05dc \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea gaussian mixture \u05e9\u05de\u05de\u05e0\u05d5 \u05e0\u05d3\u05d2\u05dd \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e2\u05e6\u05de\u05d5. \u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05dc \u05e4\u05e2\u05dd \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d4\u05e1\u05d9\u05d1\u05ea\u05d9 (\u05dc\u05d5\u05e7\u05d7 \u05d1\u05d7\u05e9\u05d1\u05d5\u05df \u05e8\u05e7 \u05d0\u05ea \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05e9\u05e0\u05d5\u05e6\u05e8\u05d5 \u05db\u05d1\u05e8) \u05d7\u05d5\u05d6\u05d4 \u05d0\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05ea\u05d5\u05d7\u05dc\u05d5\u05ea, \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05dc \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05e7\u05d5\u05d5\u05e8\u05d9\u05d0\u05e0\u05e1 \u05d0\u05dc\u05db\u05e1\u05d5\u05e0\u05d9\u05d5\u05ea \u05e9\u05dc \u05db\u05dc \u05de\u05e9\u05ea\u05e0\u05d4 \u05d1-mix \u05d5\u05de\u05e9\u05e7\u05dc\u05d9 \u05d4\u05e2\u05e8\u05d1\u05d5\u05d1). \u05d0\u05d7\u05e8\u05d9 \u05e9\u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e0\u05d7\u05d6\u05d5 \u05d5\u05e0\u05d3\u05d2\u05de\u05d5 \u05d4\u05dd \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05dc\u05d2\u05e0\u05e8\u05d5\u05d8 \u05ea\u05de\u05d5\u05e0\u05d4.\", Word.InsertLocation.after);\ntail = tail.insertParagraph(\"\u05e0\u05e6\u05d9\u05d9\u05df \u05db\u05d9 GIVT \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de-VQ-VAE \u05e0\u05d9\u05ea\u05df \u05dc\u05d0\u05de\u05df \u05d1\u05de\u05dc\u05d0\u05d5 \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05e9\u05dc\u05d8\u05e2\u05e0\u05ea \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d9\u05db\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05d1\u05e2\u05d9\u05d9\u05ea\u05d9. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd (\u05d1\u05de\u05e7\u05d5\u05dd \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d4\u05e1\u05d9\u05d1\u05ea\u05d9) \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc (\u05d4\u05e0\u05e7\u05e8\u05d0 adapter) \u05e9\u05dc Normalized Flow \u05dc\u05d2\u05e0\u05e8\u05d5\u05d8 \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2 \u05dc\u05d8\u05e0\u05d8\u05d9 \u05db\u05d5\u05dc\u05d5 \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0\u05d7\u05e8 \u05e9\u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05db\u05d1\u05e8 \u05d0\u05d5\u05de\u05e0\u05d5 \u05d5\u05db\u05db\u05d4 \u05dc\u05d4\u05e4\u05e8\u05d9\u05d3 \u05d0\u05ea \u05e9\u05e0\u05d9 \u05d4\u05e9\u05dc\u05d1\u05d9\u05dd.\", Word.InsertLocation.after);\ntail = tail.insertParagraph(\"\u05d1\u05e0\u05d5\u05e1\u05e3 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d0\u05de\u05df \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05dc\u05d0 \u05e1\u05d9\u05d1\u05ea\u05d9 \u05dc\u05d7\u05d9\u05d6\u05d5\u05d9 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e9\u05dc \u05e4\u05d0\u05e6'\u05d9\u05dd (\u05de\u05d0\u05d5\u05de\u05df \u05d3\u05d5\u05de\u05d4 \u05dcmasked languaged modeling \u05d0\u05d5 MLM). \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05d6\u05d5 \u05e9\u05dc\u05d0 \u05d9\u05d3\u05e2\u05ea\u05d9 \u05e2\u05dc\u05d9\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05d4\u05d5\u05e6\u05d2\u05d4 \u05d1\u05de\u05d0\u05de\u05e8 MaskGit.\", Word.InsertLocation.after);\ntail = tail.insertParagraph(\"https://arxiv.org/pdf/2312.02116\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Apply the review-434 rewrite: update the title/date line, rewrite the\n# body paragraphs, and append four new paragraphs (incl. the new source link)\n# that extend the review, replacing the single old arxiv link paragraph.\n\n$d = $word.ActiveDocument\n\n# Paragraph 1 - date/title line (kept as two runs joined by a manual line break)\n$d.Paragraphs(1).Range.Text = '\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 05.04.25' + [char]11 + 'GIVT: Generative Infinite-Vocabulary Transformers'\n\n# Paragraphs 2-5 - rewritten review body\n$d.Paragraphs(2).Range.Text = '\u05d4\u05d9\u05d5\u05dd \u05d7\u05d5\u05d6\u05e8\u05d9\u05dd \u05db\u05de\u05d4 \u05e9\u05e0\u05d9\u05dd \u05e9\u05e0\u05d9\u05dd \u05d0\u05d7\u05d5\u05e8\u05d4 \u05d1\u05d5 \u05de\u05d9\u05dc\u05d9\u05dd VAE, VQ-VAE, VQ-GAN \u05d4\u05d9\u05d5 \u05de\u05d5\u05e9\u05db\u05d9\u05dd \u05d0\u05d5\u05ea\u05d4 \u05ea\u05e9\u05d5\u05de\u05ea \u05dc\u05d1 \u05db\u05de\u05d5 \u05e9\u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05d4\u05d9\u05d5\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d2\u05e0\u05e8\u05d8\u05d9\u05d1\u05d9\u05d9\u05dd (\u05d0\u05de\u05e0\u05dd \u05e4\u05d7\u05d5\u05ea \u05de\u05d0\u05d2''\u05e0\u05d8\u05d9\u05dd \u05d0\u05d1\u05dc \u05d1\u05db\u05dc \u05d6\u05d0\u05ea).  \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05de\u05e6\u05d9\u05e2 \u05e9\u05db\u05dc\u05d5\u05dc \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05dc-VQ-VAE \u05e9\u05de\u05e9\u05da \u05d0\u05ea \u05e2\u05d9\u05e0\u05d9\u05d9 \u05db\u05d9 \u05db\u05d0\u05de\u05d5\u05e8 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05d1\u05e0\u05d5\u05e9\u05d0 \u05d6\u05d4 \u05d4\u05e4\u05db\u05d5 \u05dc\u05d4\u05d9\u05d5\u05ea \u05f4\u05e6\u05d9\u05e4\u05d5\u05e8 \u05e0\u05d3\u05d9\u05e8\u05d4\u05f4 \u05d1\u05e0\u05d5\u05e3 \u05e9\u05dc\u05e0\u05d5 (\u05e9\u05dc AI).'\n$d.Paragraphs(3).Range.Text = '\u05e7\u05d5\u05d3\u05dd \u05db\u05dc \u05d0\u05ea\u05df \u05d4\u05e7\u05d3\u05de\u05d4 \u05e7\u05e6\u05e8 \u05dc\u05d2\u05d1\u05d9 VQ-VAE. \u05e0\u05ea\u05d7\u05d9\u05dc \u05d0\u05ea \u05d4\u05d4\u05e1\u05d1\u05e8 \u05de-VAE \u05e9\u05d6\u05d4 \u05e8\u05d0\u05e9\u05d9 \u05ea\u05d9\u05d1\u05d5\u05ea \u05e9\u05dc Variational AutoEncoder \u05e9\u05d4\u05d5\u05de\u05e6\u05d0 \u05d0\u05d9 \u05e9\u05dd \u05d1-2014 \u05e2\u05dc \u05d9\u05d3\u05d9 Kingma \u05d4\u05d0\u05d2\u05d3\u05d9. \u05d1\u05d2\u05d3\u05d5\u05dc VAE \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05e9\u05ea\u05d9 \u05e8\u05e9\u05ea\u05d5\u05ea, \u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d3\u05e7\u05d5\u05d3\u05e8 \u05db\u05d0\u05e9\u05e8 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05de\u05e4\u05d9\u05e7 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 (\u05d0\u05d5 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2) \u05e9\u05dc \u05e4\u05d9\u05e1\u05ea \u05d3\u05d0\u05d8\u05d4 \u05db\u05d0\u05e9\u05e8 \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05dc\u05ea\u05de\u05d5\u05e0\u05d4. \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d4\u05d5\u05d0 \u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05dc \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d4\u05d2\u05d0\u05d5\u05e1\u05d9\u05ea (\u05d5\u05e7\u05d8\u05d5\u05e8 \u05ea\u05d5\u05d7\u05dc\u05d5\u05ea \u05d5\u05de\u05d8\u05e8\u05d9\u05e6\u05ea \u05e7\u05d5\u05d5\u05e8\u05d9\u05d0\u05e0\u05e1 \u05d0\u05dc\u05db\u05e1\u05d5\u05e0\u05d9\u05ea) \u05de\u05de\u05e0\u05d4 \u05d3\u05d5\u05d2\u05de\u05d9\u05dd \u05d0\u05ea \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05d4\u05de\u05d5\u05d6\u05df \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05dc\u05e9\u05d7\u05d6\u05d5\u05e8 \u05ea\u05de\u05d5\u05e0\u05ea \u05d4\u05e7\u05dc\u05d8 \u05dc\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8. '\n$d.Paragraphs(4).Range.Text = '\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05dc\u05d5\u05e1 \u05e9\u05dc VAE \u05e0\u05d1\u05e0\u05d9\u05ea \u05e2\u05dc \u05d1\u05e1\u05d9\u05e1 ELBO (\u05e9\u05d6\u05d4 Evidence Lower Bound) \u05d5\u05de\u05db\u05d9\u05dc\u05d4 2 \u05d0\u05d9\u05d1\u05e8\u05d9\u05dd. \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05dc\u05d5\u05e1 \u05d4\u05e9\u05d7\u05d6\u05d5\u05e8 \u05e9\u05dc-VAE \u05e7\u05dc\u05d0\u05e1\u05d9 \u05d4\u05d5\u05d0 \u05d4\u05e0\u05d5\u05e8\u05de\u05d4 \u05e9\u05dc \u05d4\u05e4\u05e8\u05e9 \u05d4\u05ea\u05de\u05d5\u05e0\u05d4 \u05d4\u05de\u05e9\u05d5\u05d7\u05d6\u05e8\u05ea \u05d1\u05d9\u05d7\u05e1 \u05dc\u05ea\u05de\u05d5\u05e0\u05d4 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9\u05ea (\u05d1\u05d2\u05e8\u05e1\u05d0\u05d5\u05ea \u05de\u05ea\u05e7\u05d3\u05de\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d4\u05ea\u05d5\u05d5\u05e1\u05e4\u05d5 \u05dc\u05d6\u05d4 \u05dc\u05d5\u05e1 perceptual \u05d5\u05dc\u05d5\u05e1 \u05d1\u05e1\u05d2\u05e0\u05d5\u05df GAN) \u05d5\u05d4\u05d0\u05d9\u05d1\u05e8 \u05d4\u05e9\u05e0\u05d9 \u05d4\u05d5\u05d0 KL divergence \u05d1\u05d9\u05df \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e9\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05d4\u05de\u05d5\u05e4\u05e7 \u05de\u05d4\u05d3\u05d0\u05d8\u05d4 (\u05d4\u05de\u05d9\u05d5\u05e6\u05d2 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d5\u05e7\u05d8\u05d5\u05e8 \u05ea\u05d5\u05d7\u05dc\u05d5\u05ea \u05d5\u05de\u05d8\u05e8\u05d9\u05e6\u05ea \u05e7\u05d5\u05d5\u05e8\u05d9\u05d0\u05e0\u05e1 \u05d0\u05dc\u05db\u05e1\u05d5\u05e0\u05d9\u05ea) \u05dc\u05d1\u05d9\u05df \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e0\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea. \u05d1\u05d0\u05d9\u05e0\u05e4\u05e8\u05e0\u05e1 \u05d0\u05e0\u05d5 \u05d3\u05d5\u05d2\u05de\u05d9\u05dd \u05d5\u05e7\u05d8\u05d5\u05e8 \u05dc\u05d8\u05e0\u05d8\u05d9 \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e0\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea \u05d5\u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05d5\u05ea\u05d5 \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8.'\n$d.Paragraphs(5).Range.Text = '\u05e9\u05db\u05dc\u05d5\u05dc \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e9\u05d4\u05e4\u05da \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d0\u05d5\u05d3 \u05e4\u05d5\u05e4\u05d5\u05dc\u05e8\u05d9 \u05e9\u05dc VAE \u05d4\u05d5\u05d0 VQ-VAE. \u05d1\u05de\u05e7\u05d5\u05dd \u05dc\u05d4\u05d2\u05d3\u05d9\u05e8 \u05de\u05e8\u05d7\u05d1 \u05dc\u05d8\u05e0\u05d8\u05d9 \u05d1\u05ea\u05d5\u05e8 \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d2\u05d0\u05d5\u05e1\u05d9\u05ea - \u05d0\u05dc\u05d0 \u05de\u05d2\u05d3\u05d9\u05e8 \u05d0\u05d5\u05ea\u05d5 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d9\u05e1\u05e7\u05e8\u05d8\u05d9\u05ea. \u05db\u05dc \u05e4\u05d0\u05e5'' \u05d1\u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05ea\u05d5\u05d0\u05e8(\u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9) \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d5\u05e7\u05d8\u05d5\u05e8 \u05de\u05d4-codebook \u05d1\u05d2\u05d5\u05d3\u05dc \u05e1\u05d5\u05e4\u05d9 \u05e9\u05de\u05d0\u05d5\u05de\u05df \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d4\u05d3\u05e7\u05d5\u05d3\u05e8. \u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05e9 \u05de\u05e1\u05e4\u05e8 \u05e1\u05d5\u05e4\u05d9 \u05e9\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05e4\u05d0\u05e5'' (\u05e0\u05d7\u05d6\u05d5\u05e8 \u05e2\u05dc \u05d6\u05d4 \u05e2\u05d5\u05d3 \u05de\u05e2\u05d8). \u05db\u05d0\u05e9\u05e8 \u05d0\u05d9\u05de\u05d5\u05df \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8, \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05de\u05d4-codebook \u05de\u05e1\u05ea\u05d9\u05d9\u05dd \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05e0\u05d5\u05e1\u05e3 \u05dc\u05d7\u05d9\u05d6\u05d5\u05d9 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc \u05e4\u05d0\u05e6''\u05d9\u05dd, \u05e2\u05dc \u05db\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4\u05e4\u05d0\u05e6''\u05d9\u05dd \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8. \u05de\u05d5\u05d3\u05dc \u05d6\u05d4 (\u05e0\u05d2\u05d9\u05d3 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8) \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05d7\u05d6\u05d5\u05ea \u05d1\u05d0\u05d5\u05e4\u05df \u05d0\u05d5\u05d8\u05d5\u05e8\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9 \u05d0\u05ea \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05de\u05d4-codebook (\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05e1\u05e4\u05e8\u05d5) \u05e9\u05dc \u05d4\u05e4\u05d0\u05e5'' \u05d4\u05d1\u05d0 \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05e4\u05d0\u05e6''\u05d9\u05dd \u05d4\u05e7\u05d5\u05d3\u05de\u05d9\u05dd \u05e9\u05db\u05d1\u05e8 \u05d2\u05d5\u05e0\u05e8\u05d8\u05d5. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05d4\u05d5\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e9\u05dc \u05d4\u05e4\u05d0\u05e6''\u05d9\u05dd \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05dc\u05d2\u05e0\u05e8\u05d5\u05d8 \u05d3\u05d0\u05d8\u05d4 (\u05ea\u05de\u05d5\u05e0\u05d4).'\n\n# Paragraph 6 - was the old arxiv link, now continues the rewritten review text\n$d.Paragraphs(6).Range.Text = '\u05db\u05d0\u05de\u05d5\u05e8 \u05d9\u05e9 \u05de\u05e1\u05e4\u05e8 \u05de\u05d5\u05d2\u05d1\u05dc \u05d5\u05e1\u05d5\u05e4\u05d9 \u05e9\u05dc \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05e4\u05d0\u05e5'' \u05d5\u05e9\u05d6\u05d4 \u05d3\u05d9 \u05de\u05d2\u05d1\u05d9\u05dc \u05d0\u05ea \u05d4\u05e2\u05d5\u05e9\u05e8 \u05d4\u05e1\u05de\u05e0\u05d8\u05d9 \u05e9\u05dc \u05d4\u05ea\u05de\u05d5\u05e0\u05d5\u05ea \u05e9-VQ-VAE \u05d5\u05e9\u05d9\u05d8\u05d5\u05ea \u05d3\u05d5\u05de\u05d5\u05ea \u05de\u05e1\u05d5\u05d2\u05dc\u05d5\u05ea \u05dc\u05d2\u05e0\u05e8\u05d8. \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05d4\u05de\u05e7\u05d5\u05dd \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05e2\u05db\u05e9\u05d9\u05d5 \u05de\u05d7\u05d3\u05e9 - \u05d4\u05d5\u05d0 \u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05e2\u05d1\u05d5\u05e8 \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e8\u05e6\u05d9\u05e3 (\u05dc\u05dc\u05d0 codebook) \u05e9\u05dc \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd. \u05d0\u05d1\u05dc \u05d0\u05d9\u05da \u05d0\u05e4\u05e9\u05e8 \u05dc\u05e2\u05e9\u05d5\u05ea \u05d6\u05d0\u05ea? \u05e0\u05d6\u05db\u05d9\u05e8 \u05e9-VQ-VAE \u05d0\u05e0\u05d5 \u05db\u05dc \u05e4\u05e2\u05dd \u05d7\u05d5\u05d6\u05d9\u05dd \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e7\u05d8\u05d2\u05d5\u05e8\u05d9\u05d0\u05dc\u05d9\u05ea \u05de\u05e2\u05dc \u05d4-codebook \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05d1\u05de\u05d5\u05d3\u05dc \u05d0\u05d5\u05d8\u05d5\u05e8\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9 \u05d4\u05d9\u05d0 \u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05d1\u05d2\u05d5\u05d3\u05dc \u05e9\u05dc \u05d4-codebook.'\n\n# Append four new paragraphs after the (rewritten) 6th paragraph, ending with\n# the new arxiv source link as the document's final paragraph.\n$d.Paragraphs(6).Range.InsertParagraphAfter()\n$d.Paragraphs(7).Range.Text = ' \u05d4\u05d0\u05dd \u05e0\u05d9\u05ea\u05df \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05e9\u05d9\u05d5\u05e6\u05e8 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e8\u05e6\u05d9\u05e4\u05d9\u05dd \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d0\u05d5\u05d8\u05d5\u05e8\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9\u05ea? \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05d9\u05d0 \u05db\u05df - \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05de\u05df \u05d0\u05ea \u05d4-VAE \u05d4\u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9 \u05e9\u05d4\u05e1\u05d1\u05e8\u05ea\u05d9 \u05e2\u05dc\u05d9\u05d5 \u05d1\u05ea\u05d7\u05d9\u05dc\u05ea \u05d4\u05e1\u05e7\u05d9\u05e8\u05d4. \u05d4\u05d3\u05d1\u05e8 \u05d4\u05d6\u05d4 \u05e0\u05e2\u05e9\u05d4 \u05d1\u05e8\u05de\u05d4 \u05e9\u05dc \u05e4\u05d0\u05e5'' \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc \u05ea\u05de\u05d5\u05e0\u05d4 \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05d4\u05e4\u05d0\u05e6''\u05d9\u05dd \u05e9\u05dc\u05d4. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e9\u05e0\u05d9 \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05e1\u05d9\u05d1\u05ea\u05d9 \u05e9\u05d7\u05d5\u05d6\u05d4 \u05d0\u05ea \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05d4\u05d8\u05d5\u05e7\u05df \u05d4\u05d1\u05d0 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d7\u05d9\u05d6\u05d5\u05d9 \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05dc \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea gaussian mixture \u05e9\u05de\u05de\u05e0\u05d5 \u05e0\u05d3\u05d2\u05dd \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e2\u05e6\u05de\u05d5. \u05db\u05dc\u05d5\u05de\u05e8 \u05db\u05dc \u05e4\u05e2\u05dd \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d4\u05e1\u05d9\u05d1\u05ea\u05d9 (\u05dc\u05d5\u05e7\u05d7 \u05d1\u05d7\u05e9\u05d1\u05d5\u05df \u05e8\u05e7 \u05d0\u05ea \u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05e9\u05e0\u05d5\u05e6\u05e8\u05d5 \u05db\u05d1\u05e8) \u05d7\u05d5\u05d6\u05d4 \u05d0\u05ea \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05ea\u05d5\u05d7\u05dc\u05d5\u05ea, \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e9\u05dc \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05e7\u05d5\u05d5\u05e8\u05d9\u05d0\u05e0\u05e1 \u05d0\u05dc\u05db\u05e1\u05d5\u05e0\u05d9\u05d5\u05ea \u05e9\u05dc \u05db\u05dc \u05de\u05e9\u05ea\u05e0\u05d4 \u05d1-mix \u05d5\u05de\u05e9\u05e7\u05dc\u05d9 \u05d4\u05e2\u05e8\u05d1\u05d5\u05d1). \u05d0\u05d7\u05e8\u05d9 \u05e9\u05d4\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9\u05dd \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e0\u05d7\u05d6\u05d5 \u05d5\u05e0\u05d3\u05d2\u05de\u05d5 \u05d4\u05dd \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc\u05d3\u05e7\u05d5\u05d3\u05e8 \u05dc\u05d2\u05e0\u05e8\u05d5\u05d8 \u05ea\u05de\u05d5\u05e0\u05d4.'\n\n$d.Paragraphs(7).Range.InsertParagraphAfter()\n$d.Paragraphs(8).Range.Text = '\u05e0\u05e6\u05d9\u05d9\u05df \u05db\u05d9 GIVT \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de-VQ-VAE \u05e0\u05d9\u05ea\u05df \u05dc\u05d0\u05de\u05df \u05d1\u05de\u05dc\u05d0\u05d5 \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05e9\u05dc\u05d8\u05e2\u05e0\u05ea \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d9\u05db\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05d1\u05e2\u05d9\u05d9\u05ea\u05d9. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd (\u05d1\u05de\u05e7\u05d5\u05dd \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05d4\u05e1\u05d9\u05d1\u05ea\u05d9) \u05dc\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc (\u05d4\u05e0\u05e7\u05e8\u05d0 adapter) \u05e9\u05dc Normalized Flow \u05dc\u05d2\u05e0\u05e8\u05d5\u05d8 \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2 \u05dc\u05d8\u05e0\u05d8\u05d9 \u05db\u05d5\u05dc\u05d5 \u05e9\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4 \u05dc\u05d0\u05d7\u05e8 \u05e9\u05d4\u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05d5\u05d4\u05d3\u05e7\u05d5\u05d3\u05e8 \u05db\u05d1\u05e8 \u05d0\u05d5\u05de\u05e0\u05d5 \u05d5\u05db\u05db\u05d4 \u05dc\u05d4\u05e4\u05e8\u05d9\u05d3 \u05d0\u05ea \u05e9\u05e0\u05d9 \u05d4\u05e9\u05dc\u05d1\u05d9\u05dd.'\n\n$d.Paragraphs(8).Range.InsertParagraphAfter()\n$d.Paragraphs(9).Range.Text = '\u05d1\u05e0\u05d5\u05e1\u05e3 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d0\u05de\u05df \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05dc\u05d0 \u05e1\u05d9\u05d1\u05ea\u05d9 \u05dc\u05d7\u05d9\u05d6\u05d5\u05d9 \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05dc\u05d8\u05e0\u05d8\u05d9\u05d9\u05dd \u05e9\u05dc \u05e4\u05d0\u05e6''\u05d9\u05dd (\u05de\u05d0\u05d5\u05de\u05df \u05d3\u05d5\u05de\u05d4 \u05dcmasked languaged modeling \u05d0\u05d5 MLM). \u05d4\u05d2\u05d9\u05e9\u05d4 \u05d4\u05d6\u05d5 \u05e9\u05dc\u05d0 \u05d9\u05d3\u05e2\u05ea\u05d9 \u05e2\u05dc\u05d9\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05d4\u05d5\u05e6\u05d2\u05d4 \u05d1\u05de\u05d0\u05de\u05e8 MaskGit.'\n\n$d.Paragraphs(9).Range.InsertParagraphAfter()\n$d.Paragraphs(10).Range.Text = 'https://arxiv.org/pdf/2312.02116'\n"}
